$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "percent error" rows below the existing model-results table ---
# Row 50: absolute difference between the last row (48) and the first row (17)
$ws.Range("B50").Formula = "=B48-B17"
$ws.Range("C50:G50").Formula = "=C48-C17"
$ws.Range("B50:G50").NumberFormat = "0.00"

# Row 51: percent error = row50 difference / row17 baseline
$ws.Range("B51").Formula = "=B50/B17"
$ws.Range("C51").Formula = "=C50/C17"
$ws.Range("D51").Formula = "=D50/D17"
$ws.Range("E51").Formula = "=E50/E17"
$ws.Range("F51").Formula = "=F50/F17"
$ws.Range("G51").Formula = "=G50/G17"
$ws.Range("B51:G51").NumberFormat = "0.0%"

# --- View changes: zoom to 74% and move the active selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 74
$ws.Range("K19").Select()
